$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (shifts old rows 9-14 down to 10-15),
# turning the 5-slot afternoon schedule into a 6-slot schedule with lunch
# moved from 11:30 to 12:20.
$ws.Rows.Item(9).Insert()

# Row 3: shift teaching blocks one day later (D/E/F) and clear B/C
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "MEC-3A-EAP"
$ws.Range("F3").Value = "MCT-2A-EAP"

# Row 4: Thursday class changes from MEC-3A-EAP to MCT-2A-EAP
$ws.Range("E4").Value = "MCT-2A-EAP"

# Row 6: Tuesday slot now has MEC-3A-EAP
$ws.Range("C6").Value = "MEC-3A-EAP"

# Row 8 (11:30) no longer is lunch; becomes a normal "-" slot
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"

# Row 9 (new): 12:20 lunch ("Almoço") row
$ws.Range("A9").Value = "12:20"
$ws.Range("B9").Value = "Almoço"
$ws.Range("C9").Value = "Almoço"
$ws.Range("D9").Value = "Almoço"
$ws.Range("E9").Value = "Almoço"
$ws.Range("F9").Value = "Almoço"

# Rows 10-13 keep the shifted-down times from the old rows 9-12
$ws.Range("A10").Value = "13:00"
$ws.Range("A11").Value = "13:50"
$ws.Range("A12").Value = "14:40"
$ws.Range("A13").Value = "15:30"

# Add the two new rows at the end of the table (16 and 17)
$ws.Range("A16").Value = "17:30"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"

$ws.Range("A17").Value = "18:20"

# Leave B17:F17 empty but still present in the sheet (matches the source
# file, which keeps empty <c> placeholders on the last row). Touching the
# (no-op) border setting materializes the cell without altering its style.
$ws.Range("B17:F17").Borders.LineStyle = -4142

$ws.Range("A1:F17").Columns.AutoFit() | Out-Null
